$wb = $excel.ActiveWorkbook

# --- "day" sheet: D322:D325 were stored as text (inline strings); convert
#     them to true numeric values (same digits, numeric type). ---
$day = $wb.Worksheets.Item("day")
$day.Range("D322").Value = 532541
$day.Range("D323").Value = 511218
$day.Range("D324").Value = 532478
$day.Range("D325").Value = 500469

# --- "week" sheet: append 14 new data rows (131-144) scraped from the
#     09/08/2024 11:32:24 run. ---
$week = $wb.Worksheets.Item("week")

$week.Range("A131").Value = 1
$week.Range("B131").Value = "BANKNIFTY"
$week.Range("C131").Value = "BANKNIFTY"
$week.Range("E131").Value = 0.65
$week.Range("F131").Value = 50484.5
$week.Range("G131").Value = 0
$week.Range("H131").Value = "week"
$week.Range("I131").Value = "09/08/2024 11:32:24"

$week.Range("A132").Value = 2
$week.Range("B132").Value = "DIXON"
$week.Range("C132").Value = "Dixon Technologies"
$week.Range("D132").Value = 540699
$week.Range("E132").Value = 2.5
$week.Range("F132").Value = 11740.45
$week.Range("G132").Value = 341383
$week.Range("H132").Value = "week"
$week.Range("I132").Value = "09/08/2024 11:32:24"

$week.Range("A133").Value = 3
$week.Range("B133").Value = "ULTRACEMCO"
$week.Range("C133").Value = "Ultratech Cement Limited"
$week.Range("D133").Value = 532538
$week.Range("E133").Value = 0.38
$week.Range("F133").Value = 11300.35
$week.Range("G133").Value = 189230
$week.Range("H133").Value = "week"
$week.Range("I133").Value = "09/08/2024 11:32:24"

$week.Range("A134").Value = 4
$week.Range("B134").Value = "SIEMENS"
$week.Range("C134").Value = "Siemens Limited"
$week.Range("D134").Value = 500550
$week.Range("E134").Value = 1.78
$week.Range("F134").Value = 6889.5
$week.Range("G134").Value = 725795
$week.Range("H134").Value = "week"
$week.Range("I134").Value = "09/08/2024 11:32:24"

$week.Range("A135").Value = 5
$week.Range("B135").Value = "ACC"
$week.Range("C135").Value = "Acc Limited"
$week.Range("D135").Value = 500410
$week.Range("E135").Value = -0.24
$week.Range("F135").Value = 2351.55
$week.Range("G135").Value = 191993
$week.Range("H135").Value = "week"
$week.Range("I135").Value = "09/08/2024 11:32:24"

$week.Range("A136").Value = 6
$week.Range("B136").Value = "ASTRAL"
$week.Range("C136").Value = "Astral Poly Technik Limited"
$week.Range("D136").Value = 532830
$week.Range("E136").Value = -2.44
$week.Range("F136").Value = 2031.35
$week.Range("G136").Value = 631314
$week.Range("H136").Value = "week"
$week.Range("I136").Value = "09/08/2024 11:32:24"

$week.Range("A137").Value = 7
$week.Range("B137").Value = "HDFCBANK"
$week.Range("C137").Value = "Hdfc Bank Limited"
$week.Range("D137").Value = 500180
$week.Range("E137").Value = 0.46
$week.Range("F137").Value = 1650.2
$week.Range("G137").Value = 13322309
$week.Range("H137").Value = "week"
$week.Range("I137").Value = "09/08/2024 11:32:24"

$week.Range("A138").Value = 8
$week.Range("B138").Value = "JSWSTEEL"
$week.Range("C138").Value = "Jsw Steel Limited"
$week.Range("D138").Value = 500228
$week.Range("E138").Value = 1.98
$week.Range("F138").Value = 905.1
$week.Range("G138").Value = 1722454
$week.Range("H138").Value = "week"
$week.Range("I138").Value = "09/08/2024 11:32:24"

$week.Range("A139").Value = 9
$week.Range("B139").Value = "CANFINHOME"
$week.Range("C139").Value = "Can Fin Homes Limited"
$week.Range("D139").Value = 511196
$week.Range("E139").Value = 2.26
$week.Range("F139").Value = 803.9
$week.Range("G139").Value = 663676
$week.Range("H139").Value = "week"
$week.Range("I139").Value = "09/08/2024 11:32:24"

$week.Range("A140").Value = 10
$week.Range("B140").Value = "AMBUJACEM"
$week.Range("C140").Value = "Ambuja Cements Limited"
$week.Range("D140").Value = 500425
$week.Range("E140").Value = -1.02
$week.Range("F140").Value = 632
$week.Range("G140").Value = 2606501
$week.Range("H140").Value = "week"
$week.Range("I140").Value = "09/08/2024 11:32:24"

$week.Range("A141").Value = 11
$week.Range("B141").Value = "UPL"
$week.Range("C141").Value = "Upl Limited"
$week.Range("D141").Value = 512070
$week.Range("E141").Value = 1.24
$week.Range("F141").Value = 554.6
$week.Range("G141").Value = 2061454
$week.Range("H141").Value = "week"
$week.Range("I141").Value = "09/08/2024 11:32:24"

$week.Range("A142").Value = 12
$week.Range("B142").Value = "VEDL"
$week.Range("C142").Value = "Vedanta Limited"
$week.Range("D142").Value = 500295
$week.Range("E142").Value = 1.55
$week.Range("F142").Value = 428.85
$week.Range("G142").Value = 8363641
$week.Range("H142").Value = "week"
$week.Range("I142").Value = "09/08/2024 11:32:24"

$week.Range("A143").Value = 13
$week.Range("B143").Value = "MOTHERSON"
$week.Range("C143").Value = "Motherson Sumi Systems Limited"
$week.Range("D143").Value = 517334
$week.Range("E143").Value = 2.89
$week.Range("F143").Value = 187.74
$week.Range("G143").Value = 13255204
$week.Range("H143").Value = "week"
$week.Range("I143").Value = "09/08/2024 11:32:24"

$week.Range("A144").Value = 14
$week.Range("B144").Value = "LTF"
$week.Range("C144").Value = "L&T Finance Ltd"
$week.Range("D144").Value = 533519
$week.Range("E144").Value = 0.33
$week.Range("F144").Value = 166.2
$week.Range("G144").Value = 2882574
$week.Range("H144").Value = "week"
$week.Range("I144").Value = "09/08/2024 11:32:24"
